# Ongoing abstraction of repeat lines of code:
# Insert a new parameter row (baseline_prev_labour_states) at the top of the
# parameter list (row 10), pushing every existing parameter row down by one,
# and widen column B to fit the new longer value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 10 - everything below shifts down automatically.
$ws.Rows.Item(10).Insert()

# Populate the new row with the parameter name / value pair.
$ws.Range("A10").Value = "baseline_prev_labour_states"
$ws.Range("B10").Value = "[0.04, 0.12, 0.80, 0.04]"
$ws.Range("B10").HorizontalAlignment = -4152

# Column B needs to be a bit wider to comfortably fit the new text value.
$ws.Columns.Item(2).ColumnWidth = 25.3611

# Move the active selection, matching where the editor was last working.
$ws.Range("G20").Select()
